{"js": "// Replace the 25 division-problem cells in the (only) table of the document\n// with their updated expressions. Each entry is addressed positionally by\n// (row, col) in the table grid \u2014 not by a global text search-and-replace \u2014\n// because several of the new values collide with old values that live\n// elsewhere in the table (e.g. \"58\u00f74=\" is both an original cell's text and\n// the replacement text for a different cell), so a naive find/replace could\n// double-substitute or target the wrong cell.\nconst replacements = [\n  { row: 0, col: 0, oldText: \"23\u00f74=\", newText: \"28\u00f79=\" },\n  { row: 0, col: 1, oldText: \"90\u00f73=\", newText: \"61\u00f77=\" },\n  { row: 0, col: 2, oldText: \"59\u00f75=\", newText: \"72\u00f75=\" },\n  { row: 0, col: 3, oldText: \"64\u00f72=\", newText: \"11\u00f79=\" },\n  { row: 0, col: 4, oldText: \"82\u00f72=\", newText: \"76\u00f79=\" },\n  { row: 4, col: 0, oldText: \"84\u00f75=\", newText: \"36\u00f79=\" },\n  { row: 4, col: 1, oldText: \"32\u00f75=\", newText: \"73\u00f73=\" },\n  { row: 4, col: 2, oldText: \"64\u00f77=\", newText: \"58\u00f74=\" },\n  { row: 4, col: 3, oldText: \"74\u00f78=\", newText: \"68\u00f75=\" },\n  { row: 4, col: 4, oldText: \"36\u00f78=\", newText: \"67\u00f74=\" },\n  { row: 8, col: 0, oldText: \"95\u00f76=\", newText: \"25\u00f79=\" },\n  { row: 8, col: 1, oldText: \"40\u00f75=\", newText: \"79\u00f77=\" },\n  { row: 8, col: 2, oldText: \"21\u00f74=\", newText: \"95\u00f72=\" },\n  { row: 8, col: 3, oldText: \"22\u00f76=\", newText: \"32\u00f78=\" },\n  { row: 8, col: 4, oldText: \"25\u00f73=\", newText: \"83\u00f77=\" },\n  { row: 12, col: 0, oldText: \"40\u00f79=\", newText: \"21\u00f73=\" },\n  { row: 12, col: 1, oldText: \"14\u00f77=\", newText: \"58\u00f79=\" },\n  { row: 12, col: 2, oldText: \"98\u00f77=\", newText: \"80\u00f79=\" },\n  { row: 12, col: 3, oldText: \"16\u00f75=\", newText: \"36\u00f76=\" },\n  { row: 12, col: 4, oldText: \"35\u00f78=\", newText: \"53\u00f74=\" },\n  { row: 16, col: 0, oldText: \"17\u00f72=\", newText: \"56\u00f72=\" },\n  { row: 16, col: 1, oldText: \"39\u00f79=\", newText: \"24\u00f79=\" },\n  { row: 16, col: 2, oldText: \"58\u00f74=\", newText: \"28\u00f75=\" },\n  { row: 16, col: 3, oldText: \"90\u00f72=\", newText: \"15\u00f74=\" },\n  { row: 16, col: 4, oldText: \"94\u00f78=\", newText: \"45\u00f72=\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst cells = replacements.map((r) => table.getCell(r.row, r.col));\ncells.forEach((cell) => cell.load(\"value\"));\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const { oldText, newText } = replacements[i];\n  const cell = cells[i];\n  // Defensive check: confirm we're editing the cell we think we are before\n  // overwriting it (should always hold true for this document).\n  if (cell.value !== oldText) {\n    console.log(\n      `Warning: cell (${replacements[i].row}, ${replacements[i].col}) was \"${cell.value}\", expected \"${oldText}\"`\n    );\n  }\n  cell.value = newText;\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 division-problem cells in the (only) table of the document\n# with their updated expressions. Each entry is addressed positionally by\n# (row, col) in the table grid (1-based, per Word COM conventions) -- not by\n# a global text search-and-replace -- because several of the new values\n# collide with old values that live elsewhere in the table (e.g. \"58\u00f74=\" is\n# both an original cell's text and the replacement text for a different\n# cell), so a naive find/replace could double-substitute or target the\n# wrong cell.\n$replacements = @(\n    @{ Row = 1; Col = 1; OldText = \"23\u00f74=\"; NewText = \"28\u00f79=\" },\n    @{ Row = 1; Col = 2; OldText = \"90\u00f73=\"; NewText = \"61\u00f77=\" },\n    @{ Row = 1; Col = 3; OldText = \"59\u00f75=\"; NewText = \"72\u00f75=\" },\n    @{ Row = 1; Col = 4; OldText = \"64\u00f72=\"; NewText = \"11\u00f79=\" },\n    @{ Row = 1; Col = 5; OldText = \"82\u00f72=\"; NewText = \"76\u00f79=\" },\n    @{ Row = 5; Col = 1; OldText = \"84\u00f75=\"; NewText = \"36\u00f79=\" },\n    @{ Row = 5; Col = 2; OldText = \"32\u00f75=\"; NewText = \"73\u00f73=\" },\n    @{ Row = 5; Col = 3; OldText = \"64\u00f77=\"; NewText = \"58\u00f74=\" },\n    @{ Row = 5; Col = 4; OldText = \"74\u00f78=\"; NewText = \"68\u00f75=\" },\n    @{ Row = 5; Col = 5; OldText = \"36\u00f78=\"; NewText = \"67\u00f74=\" },\n    @{ Row = 9; Col = 1; OldText = \"95\u00f76=\"; NewText = \"25\u00f79=\" },\n    @{ Row = 9; Col = 2; OldText = \"40\u00f75=\"; NewText = \"79\u00f77=\" },\n    @{ Row = 9; Col = 3; OldText = \"21\u00f74=\"; NewText = \"95\u00f72=\" },\n    @{ Row = 9; Col = 4; OldText = \"22\u00f76=\"; NewText = \"32\u00f78=\" },\n    @{ Row = 9; Col = 5; OldText = \"25\u00f73=\"; NewText = \"83\u00f77=\" },\n    @{ Row = 13; Col = 1; OldText = \"40\u00f79=\"; NewText = \"21\u00f73=\" },\n    @{ Row = 13; Col = 2; OldText = \"14\u00f77=\"; NewText = \"58\u00f79=\" },\n    @{ Row = 13; Col = 3; OldText = \"98\u00f77=\"; NewText = \"80\u00f79=\" },\n    @{ Row = 13; Col = 4; OldText = \"16\u00f75=\"; NewText = \"36\u00f76=\" },\n    @{ Row = 13; Col = 5; OldText = \"35\u00f78=\"; NewText = \"53\u00f74=\" },\n    @{ Row = 17; Col = 1; OldText = \"17\u00f72=\"; NewText = \"56\u00f72=\" },\n    @{ Row = 17; Col = 2; OldText = \"39\u00f79=\"; NewText = \"24\u00f79=\" },\n    @{ Row = 17; Col = 3; OldText = \"58\u00f74=\"; NewText = \"28\u00f75=\" },\n    @{ Row = 17; Col = 4; OldText = \"90\u00f72=\"; NewText = \"15\u00f74=\" },\n    @{ Row = 17; Col = 5; OldText = \"94\u00f78=\"; NewText = \"45\u00f72=\" }\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\nforeach ($r in $replacements) {\n    $cell = $tbl.Cell($r.Row, $r.Col)\n    # Cell.Range.Text includes a trailing cell-mark (CR + BEL); strip it off\n    # before comparing against the expected original text.\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $r.OldText) {\n        Write-Output \"Warning: cell ($($r.Row), $($r.Col)) was '$current', expected '$($r.OldText)'\"\n    }\n    $cell.Range.Text = $r.NewText\n}\n"}
